$d = $word.ActiveDocument

function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq $text) {
            return $i
        }
    }
    return -1
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParaXml($index, $bodyXml) {
    $p = $d.Paragraphs.Item($index)
    $xml = $pkgHeader + '<w:body>' + $bodyXml + '</w:body>' + $pkgFooter
    $p.Range.InsertXML($xml)
}

function Insert-ParaXmlBefore($index, $bodyXml) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertParagraphBefore() | Out-Null
    $newP = $d.Paragraphs.Item($index)
    $xml = $pkgHeader + '<w:body>' + $bodyXml + '</w:body>' + $pkgFooter
    $newP.Range.InsertXML($xml)
}

# 1. "TODO: MAP" -> "{{map:priorities}}" with reworked paragraph/run formatting
$i = Find-ParaIndex "TODO: MAP"
Set-ParaXml $i '<w:p w14:paraId="2A7AE73F" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Normal1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:noProof/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00363ECD"><w:rPr><w:noProof/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>{{map:priorities}}</w:t></w:r></w:p>'

# 2. Remove lastRenderedPageBreak from "Map of blueprint priorities..." paragraph
$i = Find-ParaIndex "Map of blueprint priorities in {{value:summary_unit_name}}"
Set-ParaXml $i '<w:p w14:paraId="48062B29" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Normal1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:t>Map of blueprint priorities in {{value:summary_unit_name}}</w:t></w:r></w:p>'

# 3. Add lastRenderedPageBreak to "Priority categories" paragraph
$i = Find-ParaIndex "Priority categories"
Set-ParaXml $i '<w:p w14:paraId="2FB2A47A" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_umkkwt4fqxa" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Priority categories</w:t></w:r></w:p>'

# 4. Insert a page-break paragraph (carrying bookmark _qisommxqtwh2, id 2) before "Indicators",
#    then strip that bookmark from "Indicators" itself and add lastRenderedPageBreak there.
$i = Find-ParaIndex "Indicators"
Insert-ParaXmlBefore $i '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="2" w:name="_qisommxqtwh2" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="2"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:br w:type="page"/></w:r></w:p>'
$i = Find-ParaIndex "Indicators"
Set-ParaXml $i '<w:p w14:paraId="7837CA76" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:lastRenderedPageBreak/><w:t>Indicators</w:t></w:r></w:p>'

# 5. Remove the stray _GoBack bookmark from the ecosystems table caption (it gets relocated to doc end)
$i = Find-ParaIndex "{{table:ecosystems}}"
Set-ParaXml $i '<w:p w14:paraId="54BB0A7B" w14:textId="22BA6E5F" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00FA5ABE"><w:pPr><w:pStyle w:val="TableCaption"/></w:pPr><w:r><w:t>{{table:ecosystems}}</w:t></w:r></w:p>'

# 6. Insert a page-break paragraph (carrying bookmark _7hn4qmarsl55, renumbered id 4) before "Threats",
#    then strip the old bookmark (old id 5) from "Threats" itself (lastRenderedPageBreak already present there).
$i = Find-ParaIndex "Threats"
Insert-ParaXmlBefore $i '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="4" w:name="_7hn4qmarsl55" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="4"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:br w:type="page"/></w:r></w:p>'
$i = Find-ParaIndex "Threats"
Set-ParaXml $i '<w:p w14:paraId="2CB2F0A1" w14:textId="5BFF8546" w:rsidR="00673AED" w:rsidRPr="00C01443" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/></w:rPr></w:pPr><w:r w:rsidRPr="00C01443"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:lastRenderedPageBreak/><w:t>Threats</w:t></w:r></w:p>'

# 7. Renumber the remaining bookmark ids down by one (6->5, 7->6, 8->7, 9->8, 10->9, 11->10, 12->11)
$i = Find-ParaIndex "Sea level rise"
Set-ParaXml $i '<w:p w14:paraId="67C62515" w14:textId="77777777" w:rsidR="00673AED" w:rsidRPr="00C01443" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="5" w:name="_vfv3vhs1u8zh" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="5"/><w:r w:rsidRPr="00C01443"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Sea level rise</w:t></w:r></w:p>'

$i = Find-ParaIndex "Urban growth"
Set-ParaXml $i '<w:p w14:paraId="1A7BC10D" w14:textId="77777777" w:rsidR="00673AED" w:rsidRPr="00C01443" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="6" w:name="_vv7otamh85pf" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="6"/><w:r w:rsidRPr="00C01443"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Urban growth</w:t></w:r></w:p>'

$i = Find-ParaIndex "Partners"
Set-ParaXml $i '<w:p w14:paraId="44088952" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/></w:rPr></w:pPr><w:bookmarkStart w:id="7" w:name="_ca20d8dihywo" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="7"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:t>Partners</w:t></w:r></w:p>'

$i = Find-ParaIndex "{{PARTNERS}}"
Set-ParaXml $i '<w:p w14:paraId="1688DFA4" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Normal1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="8" w:name="_a1ytnl6sgu9x" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="8"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{PARTNERS}}</w:t></w:r></w:p>'

$i = Find-ParaIndex "Ownership"
Set-ParaXml $i '<w:p w14:paraId="1BF849AF" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading1"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/></w:rPr></w:pPr><w:bookmarkStart w:id="9" w:name="_51kemobqovsy" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="9"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/></w:rPr><w:t>Ownership</w:t></w:r></w:p>'

$i = Find-ParaIndex "Conserved lands ownership"
Set-ParaXml $i '<w:p w14:paraId="1F2F6D68" w14:textId="77777777" w:rsidR="00673AED" w:rsidRDefault="00673AED" w:rsidP="00673AED"><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="10" w:name="_6o8u7emblwbs" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="10"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Conserved lands ownership</w:t></w:r></w:p>'

$i = Find-ParaIndex "Land protection status"
Set-ParaXml $i '<w:p w14:paraId="3528E969" w14:textId="7FDED608" w:rsidR="00673AED" w:rsidRPr="00C01443" w:rsidRDefault="00673AED" w:rsidP="00C01443"><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="9" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="312" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="11" w:name="_il78auds0hi7" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="11"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Land protection status</w:t></w:r></w:p>'

# 8. Re-add the _GoBack bookmark at the very end of the document (right after {{table:protection}})
$i = Find-ParaIndex "{{table:protection}}"
$p = $d.Paragraphs.Item($i)
$endPoint = $d.Range($p.Range.End, $p.Range.End)
$d.Bookmarks.Add("_GoBack", $endPoint) | Out-Null

Write-Output "done"
